$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1758430.6
$ws.Range("I15").Value = 1758430.6
$ws.Range("K15").Value = 5275291.800000001
$ws.Range("M15").Value = -5275122.800000001
$ws.Range("H28").Value = 791
$ws.Range("I28").Value = 674.375
$ws.Range("J28").Value = 1102
$ws.Range("K28").Value = 674.375
$ws.Range("L28").Value = 1102
$ws.Range("M28").Value = -189.375
$ws.Range("N28").Value = -2072
$ws.Range("H62").Value = 38108.1
$ws.Range("I62").Value = 8260
$ws.Range("J62").Value = 67956.2
$ws.Range("K62").Value = 8260
$ws.Range("L62").Value = 67956.2
$ws.Range("M62").Value = -7636
$ws.Range("N62").Value = -69204.2
$ws.Range("H65").Value = 38108.1
$ws.Range("I65").Value = 8260
$ws.Range("J65").Value = 67956.2
$ws.Range("K65").Value = 41300
$ws.Range("L65").Value = 339781
$ws.Range("M65").Value = -38180
$ws.Range("N65").Value = -346021
$ws.Range("H94").Value = 2621.1667
$ws.Range("I94").Value = 1431.75
$ws.Range("K94").Value = 1431.75
$ws.Range("M94").Value = -980.75
$ws.Range("H111").Value = 989.8570999999999
$ws.Range("I111").Value = 989.8570999999999
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2969.5713
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = 97.42870000000039
$ws.Range("H113").Value = 5482.0415
$ws.Range("I113").Value = 5599.3887
$ws.Range("J113").Value = 5130
$ws.Range("K113").Value = 5599.3887
$ws.Range("L113").Value = 5130
$ws.Range("M113").Value = -2345.3887
$ws.Range("N113").Value = -11638
$ws.Range("H135").Value = 845.9231
$ws.Range("I135").Value = 419.5
$ws.Range("J135").Value = 2267.3333
$ws.Range("K135").Value = 3775.5
$ws.Range("L135").Value = 20405.9997
$ws.Range("M135").Value = -1240.5
$ws.Range("N135").Value = -25475.9997
$ws.Range("H137").Value = 2058.9395
$ws.Range("I137").Value = 2409.1765
$ws.Range("J137").Value = 1686.8125
$ws.Range("K137").Value = 7227.529500000001
$ws.Range("L137").Value = 5060.4375
$ws.Range("M137").Value = -4677.529500000001
$ws.Range("N137").Value = -10160.4375
$ws.Range("H138").Value = 2206.946
$ws.Range("I138").Value = 1138.1613
$ws.Range("J138").Value = 7729
$ws.Range("K138").Value = 3414.4839
$ws.Range("L138").Value = 23187
$ws.Range("M138").Value = 1725.5161
$ws.Range("N138").Value = -33467
$ws.Range("H140").Value = 51623.332
$ws.Range("J140").Value = 51623.332
$ws.Range("L140").Value = 51623.332
$ws.Range("N140").Value = -61983.332

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4640.8203
$ws.Range("I32").Value = 4156.0615
$ws.Range("J32").Value = 9549
$ws.Range("K32").Value = 4156.0615
$ws.Range("L32").Value = 9549
$ws.Range("M32").Value = -3869.0615
$ws.Range("N32").Value = -10123
$ws.Range("H45").Value = 3690.5264
$ws.Range("I45").Value = 1448.5
$ws.Range("J45").Value = 5321.091
$ws.Range("K45").Value = 1448.5
$ws.Range("L45").Value = 5321.091
$ws.Range("M45").Value = -1071.5
$ws.Range("N45").Value = -6075.091
$ws.Range("H61").Value = 2573.8635
$ws.Range("I61").Value = 1850.7778
$ws.Range("K61").Value = 1850.7778
$ws.Range("M61").Value = -1638.7778
$ws.Range("H123").Value = 44500
$ws.Range("J123").Value = 44500
$ws.Range("L123").Value = 44500
$ws.Range("N123").Value = -54300
$ws.Range("H136").Value = 2573.8635
$ws.Range("I136").Value = 1850.7778
$ws.Range("K136").Value = 5552.3334
$ws.Range("M136").Value = -3002.3334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4659.077
$ws.Range("I134").Value = 2286.8
$ws.Range("J134").Value = 6855.6294
$ws.Range("K134").Value = 6860.400000000001
$ws.Range("L134").Value = 20566.8882
$ws.Range("M134").Value = -4325.400000000001
$ws.Range("N134").Value = -25636.8882

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8774403
$ws.Range("I31").Value = 1713.2858
$ws.Range("J31").Value = 19611256
$ws.Range("K31").Value = 1713.2858
$ws.Range("L31").Value = 19611256
$ws.Range("M31").Value = -1418.2858
$ws.Range("N31").Value = -19611846
$ws.Range("H34").Value = 8774403
$ws.Range("I34").Value = 1713.2858
$ws.Range("J34").Value = 19611256
$ws.Range("K34").Value = 1713.2858
$ws.Range("L34").Value = 19611256
$ws.Range("M34").Value = -1511.2858
$ws.Range("N34").Value = -19611660
$ws.Range("H122").Value = 1572.4814
$ws.Range("I122").Value = 1057.591
$ws.Range("J122").Value = 3838
$ws.Range("K122").Value = 3172.773
$ws.Range("L122").Value = 11514
$ws.Range("M122").Value = -722.7729999999997
$ws.Range("N122").Value = -16414

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 166.61111
$ws.Range("I12").Value = 70.333336
$ws.Range("J12").Value = 185.86667
$ws.Range("K12").Value = 211.000008
$ws.Range("L12").Value = 557.60001
$ws.Range("M12").Value = -38.00000800000001
$ws.Range("N12").Value = -903.60001
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1500
$ws.Range("K16").Value = 4500
$ws.Range("M16").Value = -4327
$ws.Range("H34").Value = 11009.333
$ws.Range("J34").Value = 11009.333
$ws.Range("L34").Value = 33027.999
$ws.Range("N34").Value = -33195.999
$ws.Range("H39").Value = 1326.5333
$ws.Range("J39").Value = 1207
$ws.Range("L39").Value = 3621
$ws.Range("N39").Value = -4209
$ws.Range("H55").Value = 15808.143
$ws.Range("J55").Value = 18277.834
$ws.Range("L55").Value = 54833.50199999999
$ws.Range("N55").Value = -55187.50199999999
$ws.Range("H106").Value = 3369.2307
$ws.Range("J106").Value = 3369.2307
$ws.Range("L106").Value = 10107.6921
$ws.Range("N106").Value = -11999.6921
$ws.Range("H131").Value = 1134.2106
$ws.Range("J131").Value = 1134.2106
$ws.Range("L131").Value = 3402.6318
$ws.Range("N131").Value = -13482.6318

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1669606
$ws.Range("I132").Value = 2978389.5
$ws.Range("J132").Value = 3881.6365
$ws.Range("K132").Value = 8935168.5
$ws.Range("L132").Value = 11644.9095
$ws.Range("M132").Value = -8932638.5
$ws.Range("N132").Value = -16704.9095
$ws.Range("H141").Value = 28980
$ws.Range("J141").Value = 28980
$ws.Range("L141").Value = 28980
$ws.Range("N141").Value = -39340

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 47591.668
$ws.Range("J134").Value = 47591.668
$ws.Range("L134").Value = 47591.668
$ws.Range("N134").Value = -57731.668

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4701.4116
$ws.Range("I122").Value = 3646.0454
$ws.Range("J122").Value = 6636.25
$ws.Range("K122").Value = 10938.1362
$ws.Range("L122").Value = 19908.75
$ws.Range("M122").Value = -8488.136200000001
$ws.Range("N122").Value = -24808.75
$ws.Range("H140").Value = 35371.6
$ws.Range("J140").Value = 35371.6
$ws.Range("L140").Value = 35371.6
$ws.Range("N140").Value = -45731.6
$ws.Range("H141").Value = 35016.43
$ws.Range("J141").Value = 35016.43
$ws.Range("L141").Value = 35016.43
$ws.Range("N141").Value = -45376.43
